$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting for Price cells whose new values would otherwise be
# auto-converted to numbers by Excel (so they match the original text storage).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"

# Apply updated values
$ws.Range("D2").Value = '67.096.56'
$ws.Range("E2").Value = '  -0.35%  '
$ws.Range("D3").Value = '2.616.05'
$ws.Range("E3").Value = '  -0.97%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '589.56'
$ws.Range("E5").Value = '  -1.48%  '
$ws.Range("D6").Value = '165.55'
$ws.Range("E6").Value = '  -0.59%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '0.532'
$ws.Range("D9").Value = '2.616.50'
$ws.Range("E9").Value = '  -0.94%  '
$ws.Range("E10").Value = '  -4.15%  '
$ws.Range("D11").Value = '0.160'
$ws.Range("E11").Value = '  +1.51%  '
$ws.Range("E12").Value = '  -0.24%  '
$ws.Range("D13").Value = '5.20'
$ws.Range("E13").Value = '  -0.50%  '
$ws.Range("D14").Value = '27.29'
$ws.Range("E14").Value = '  -2.50%  '
$ws.Range("D15").Value = '3.092.21'
$ws.Range("E15").Value = '  -1.03%  '
$ws.Range("D16").Value = '0.0000180'
$ws.Range("E16").Value = '  -2.34%  '
$ws.Range("D17").Value = '67.073.00'
$ws.Range("E17").Value = '  -0.46%  '
$ws.Range("D18").Value = '2.632.60'
$ws.Range("E18").Value = '  -0.45%  '
$ws.Range("D19").Value = '11.77'
$ws.Range("E19").Value = '  -1.04%  '
$ws.Range("D20").Value = '7.81'
$ws.Range("E20").Value = '  -1.08%  '
$ws.Range("D21").Value = '354.83'
$ws.Range("E21").Value = '  -2.47%  '
$ws.Range("E22").Value = '  -2.80%  '
$ws.Range("E23").Value = '  -2.96%  '
$ws.Range("D24").Value = '10.51'
$ws.Range("E24").Value = '  -5.32%  '
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("E26").Value = '  -4.72%  '
$ws.Range("D27").Value = '69.09'
$ws.Range("E27").Value = '  -2.68%  '
$ws.Range("D28").Value = '2.746.61'
$ws.Range("E28").Value = '  -1.09%  '
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.00%  '
$ws.Range("D30").Value = '0.0₃0995'
$ws.Range("E30").Value = '  -2.73%  '
$ws.Range("D31").Value = '543.13'
$ws.Range("E32").Value = '  -2.30%  '
$ws.Range("E33").Value = '  -3.31%  '
$ws.Range("E34").Value = '  -2.71%  '
$ws.Range("E35").Value = '  +1.27%  '
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  +0.01%  '
$ws.Range("E37").Value = '  -3.66%  '
$ws.Range("D38").Value = '157.12'
$ws.Range("E38").Value = '  -0.52%  '
$ws.Range("E39").Value = '  -2.50%  '
$ws.Range("E40").Value = '  -2.25%  '
$ws.Range("E41").Value = '  +1.79%  '
$ws.Range("E42").Value = '  -1.08%  '
$ws.Range("E43").Value = '  -2.38%  '
$ws.Range("E45").Value = '  -4.10%  '
$ws.Range("E46").Value = '  -1.05%  '
$ws.Range("D47").Value = '151.34'
$ws.Range("E47").Value = '  -1.89%  '
$ws.Range("E48").Value = '  -3.36%  '
$ws.Range("E49").Value = '  -2.89%  '
$ws.Range("E50").Value = '  -1.69%  '
$ws.Range("E51").Value = '  -1.13%  '

Write-Host "Updated cryptos list"